$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$i0 = @(8, 8, 8, 5, 8, 7, 8, 6, 6, 8, 6, 5, 7, 8, 8, 6)
$if = @(8, 8, 8, 5, 8, 7, 8, 6, 7, 9, 6, 6, 7, 8, 8, 6)

for ($r = 0; $r -lt 16; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
